$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Row 2 (CHROME) now points at the public PPM demo URL.
$ws.Range("B2").Value = "http://ppmdemo.mfadvantageinc.com/menu.html"

# New row 3: MSEDGE pointing at the (shortened) nimbusserver URL, copying
# the same cell formatting used by row 2 so borders/styles match.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("A3").Value = "MSEDGE"
$ws.Range("B3").Value = "http://nimbusserver.aos.com:8088"

# Column B needs to widen to fit the longer text now in it.
$ws.Columns.Item(2).ColumnWidth = 41
